$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 14) to the Storage table
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = "SM_FUR_VS_AB"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 1000
$ws.Range("H14").Value = 4010
$ws.Range("I14").Value = 1009

# Match the formatting (vertical-center style) used by the other rows in columns H and I
$ws.Range("H14").VerticalAlignment = -4108
$ws.Range("I14").VerticalAlignment = -4108

# Move/update the active selection as it appears after entering the new row
$ws.Range("J15").Select()
